$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, copying the style from the existing
# header cell H1 so the bold/border/centered formatting carries over.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Values for the new I/J columns, row by row (rows 2-17).
$iValues = @(1, 4, 7, 5, 1, 4, 1, 1, 4, 1, 1, 1, 5, 3, 3, 1)
$jValues = @(2, 5, 8, 8, 5, 7, 4, 4, 6, 4, 4, 3, 8, 5, 4, 2)

for ($r = 0; $r -lt 16; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
